$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("P18").Value = 187.29
$wsVentasGrupo.Range("Q18").Value = 44.6
$wsVentasGrupo.Range("Q32").Value = "1 de 30"

# --- Sheet: VENTA MENSUAL ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F18").Value = 1743.32
$wsVentaMensual.Range("F32").Value = 14489.7

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column F width: stored sheet width 24 -> 25 (Excel ColumnWidth is offset
# from the stored width by 5/6 of a character due to internal cell padding).
$wsCumplimiento.Columns.Item(6).ColumnWidth = 25 - (5/6)

$wsCumplimiento.Range("D10").Value = 228.89
$wsCumplimiento.Range("E10").Value = 421.36
$wsCumplimiento.Range("F10").Value = 0.3520030757400999

$wsCumplimiento.Range("D14").Value = 44.6
$wsCumplimiento.Range("E14").Value = 438.4
$wsCumplimiento.Range("F14").Value = 0.09233954451345756

$wsCumplimiento.Range("D19").Value = 14483.94
$wsCumplimiento.Range("E19").Value = 15053.85107555787
$wsCumplimiento.Range("F19").Value = 0.4903528487607615
